# "Generate Report for Handback" — refresh the localization-status report
# after a handback cycle completes:
#   * Overview sheet: status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   * zh-cn sheet: the handback timestamp advances
#   * de-de sheet: the handback filename/timestamp columns advance and the
#     stale "version mismatch" error on row 3 is cleared
#   * a handful of columns get re-sized (as if Excel had just reflowed the
#     report after the data refresh)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Cells.Item(2, 5).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(2, 6).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(3, 5).Value = "Handed back: in sync with en-US"
$overview.Cells.Item(3, 6).Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Latest Handback DateTime (column K) refreshed for both rows
$zhcn.Cells.Item(2, 11).Value = "2016-09-07 10:16:27"
$zhcn.Cells.Item(3, 11).Value = "2016-09-07 10:16:27"

# Error Detail (column P) on row 3 is now clear - the handback is current
$zhcn.Cells.Item(3, 16).Value = ""

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Latest Handback DateTime (column K) refreshed for both rows
$dede.Cells.Item(2, 11).Value = "2016-09-07 10:16:44"
$dede.Cells.Item(3, 11).Value = "2016-09-07 10:16:44"

# Error Detail (column P) on row 3 is now clear - the handback is current
$dede.Cells.Item(3, 16).Value = ""

# ---------------------------------------------------------------------
# Column width refresh (Excel quantizes ColumnWidth to 1/6-character
# steps, so pick the input that lands on the nearest reachable width)
# ---------------------------------------------------------------------
function Set-ClosestColumnWidth($column, $targetWidth) {
    $nearestSixth = [Math]::Round($targetWidth * 6.0) / 6.0
    $column.ColumnWidth = $nearestSixth - (5.0 / 6.0)
}

Set-ClosestColumnWidth $overview.Columns.Item(5) 29.9777050018311
Set-ClosestColumnWidth $overview.Columns.Item(6) 29.9777050018311

Set-ClosestColumnWidth $zhcn.Columns.Item(3) 29.9777050018311
Set-ClosestColumnWidth $zhcn.Columns.Item(16) 13.7470531463623

Set-ClosestColumnWidth $dede.Columns.Item(3) 29.9777050018311
Set-ClosestColumnWidth $dede.Columns.Item(16) 13.7470531463623
